$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Slit2"
$ws.Range("C2").Value = "Robo1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.01246433333333333
$ws.Range("H2").Value = 0.037393
$ws.Range("I2").Value = 0.0065371131913745
$ws.Range("J2").Value = 0.006537113191374499
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1325553333333333
$ws.Range("N2").Value = 0.397666
$ws.Range("O2").Value = 0.00533964316398423
$ws.Range("P2").Value = 0.00533964316398423
$ws.Range("Q2").Value = 0.001652213859777778
$ws.Range("R2").Value = 0.014869924738
$ws.Range("S2").Value = 0.00003490585176451398
$ws.Range("T2").Value = 0.00003490585176451397

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Slit2"
$ws.Range("C3").Value = "Robo1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.01246433333333333
$ws.Range("H3").Value = 0.037393
$ws.Range("I3").Value = 0.0065371131913745
$ws.Range("J3").Value = 0.006537113191374499
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 17.178266
$ws.Range("N3").Value = 51.534798
$ws.Range("O3").Value = 0.6919812904497951
$ws.Range("P3").Value = 0.691981290449795
$ws.Range("Q3").Value = 0.2141156335126667
$ws.Range("R3").Value = 1.927040701614
$ws.Range("S3").Value = 0.004523560021983705
$ws.Range("T3").Value = 0.004523560021983703

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Slit2"
$ws.Range("C4").Value = "Robo1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.01246433333333333
$ws.Range("H4").Value = 0.037393
$ws.Range("I4").Value = 0.0065371131913745
$ws.Range("J4").Value = 0.006537113191374499
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 7.513933666666667
$ws.Range("N4").Value = 22.541801
$ws.Range("O4").Value = 0.3026790663862208
$ws.Range("P4").Value = 0.3026790663862208
$ws.Range("Q4").Value = 0.09365617386588888
$ws.Range("R4").Value = 0.8429055647929999
$ws.Range("S4").Value = 0.001978647317626282
$ws.Range("T4").Value = 0.001978647317626282

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Slit2"
$ws.Range("C5").Value = "Robo1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.9943730000000001
$ws.Range("H5").Value = 2.983119
$ws.Range("I5").Value = 0.5215143627507798
$ws.Range("J5").Value = 0.5215143627507798
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1325553333333333
$ws.Range("N5").Value = 0.397666
$ws.Range("O5").Value = 0.00533964316398423
$ws.Range("P5").Value = 0.00533964316398423
$ws.Range("Q5").Value = 0.1318094444726667
$ws.Range("R5").Value = 1.186285000254
$ws.Range("S5").Value = 0.002784700601981793
$ws.Range("T5").Value = 0.002784700601981793

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Slit2"
$ws.Range("C6").Value = "Robo1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.9943730000000001
$ws.Range("H6").Value = 2.983119
$ws.Range("I6").Value = 0.5215143627507798
$ws.Range("J6").Value = 0.5215143627507798
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 17.178266
$ws.Range("N6").Value = 51.534798
$ws.Range("O6").Value = 0.6919812904497951
$ws.Range("P6").Value = 0.691981290449795
$ws.Range("Q6").Value = 17.081603897218
$ws.Range("R6").Value = 153.734435074962
$ws.Range("S6").Value = 0.3608781817243871
$ws.Range("T6").Value = 0.3608781817243871

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Slit2"
$ws.Range("C7").Value = "Robo1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.9943730000000001
$ws.Range("H7").Value = 2.983119
$ws.Range("I7").Value = 0.5215143627507798
$ws.Range("J7").Value = 0.5215143627507798
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 7.513933666666667
$ws.Range("N7").Value = 22.541801
$ws.Range("O7").Value = 0.3026790663862208
$ws.Range("P7").Value = 0.3026790663862208
$ws.Range("Q7").Value = 7.471652761924334
$ws.Range("R7").Value = 67.244874857319
$ws.Range("S7").Value = 0.1578514804244109
$ws.Range("T7").Value = 0.1578514804244109

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Slit2"
$ws.Range("C8").Value = "Robo1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.8998656666666666
$ws.Range("H8").Value = 2.699597
$ws.Range("I8").Value = 0.4719485240578458
$ws.Range("J8").Value = 0.4719485240578457
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1325553333333333
$ws.Range("N8").Value = 0.397666
$ws.Range("O8").Value = 0.00533964316398423
$ws.Range("P8").Value = 0.00533964316398423
$ws.Range("Q8").Value = 0.1192819934002222
$ws.Range("R8").Value = 1.073537940602
$ws.Range("S8").Value = 0.002520036710237923
$ws.Range("T8").Value = 0.002520036710237923

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Slit2"
$ws.Range("C9").Value = "Robo1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.8998656666666666
$ws.Range("H9").Value = 2.699597
$ws.Range("I9").Value = 0.4719485240578458
$ws.Range("J9").Value = 0.4719485240578457
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 17.178266
$ws.Range("N9").Value = 51.534798
$ws.Range("O9").Value = 0.6919812904497951
$ws.Range("P9").Value = 0.691981290449795
$ws.Range("Q9").Value = 15.45813178626733
$ws.Range("R9").Value = 139.123186076406
$ws.Range("S9").Value = 0.3265795487034243
$ws.Range("T9").Value = 0.3265795487034241

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Slit2"
$ws.Range("C10").Value = "Robo1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.8998656666666666
$ws.Range("H10").Value = 2.699597
$ws.Range("I10").Value = 0.4719485240578458
$ws.Range("J10").Value = 0.4719485240578457
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 7.513933666666667
$ws.Range("N10").Value = 22.541801
$ws.Range("O10").Value = 0.3026790663862208
$ws.Range("P10").Value = 0.3026790663862208
$ws.Range("Q10").Value = 6.761530928244111
$ws.Range("R10").Value = 60.853778354197
$ws.Range("S10").Value = 0.1428489386441836
$ws.Range("T10").Value = 0.1428489386441836
